$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HBAN")

# Widen column B to match columns C/D (16.5).
# NOTE: Excel's ColumnWidth setter quantizes to whole pixels (stored width =
# (round(CharWidth*7)+5)/7 for the default Arial 11/MDW=7 font), so the
# stored width can only land on multiples of 1/7 plus that padding; 16.5
# itself (…*7=115.5) is not reachable exactly. 15.86 is the closest input
# that rounds to the nearest achievable stored width (~16.571428571428573).
$ws.Columns.Item(2).ColumnWidth = 15.86

# Fill in previously-blank values in column B (2020 period) across the cash-flow rows
$ws.Range("B3").Value = 319000000.0
$ws.Range("B4").Value = 642000000.0
$ws.Range("B5").Value = -920000000.0
$ws.Range("B6").Value = -15000000.0
$ws.Range("B8").Value = -611000000.0
$ws.Range("B10").Value = 1254000000.0
$ws.Range("B11").Value = -130000000.0
$ws.Range("B13").Value = 134000000.0
$ws.Range("B14").Value = -2629000000.0
$ws.Range("B15").Value = 1521000000.0
$ws.Range("B16").Value = -5937000000.0
$ws.Range("B17").Value = -4586000000.0
$ws.Range("B19").Value = 1473000000.0
$ws.Range("B20").Value = -713000000.0
$ws.Range("B21").Value = 4680000000.0
$ws.Range("B22").Value = 11505000000.0
$ws.Range("B23").Value = 6822000000.0
$ws.Range("B24").Value = 6595000000.0
$ws.Range("B25").Value = 13417000000.0
$ws.Range("B26").Value = 90000000.0
$ws.Range("B28").Value = -1546000000.0
$ws.Range("B29").Value = -2495000000.0
$ws.Range("B30").Value = -101000000.0
$ws.Range("B31").Value = 1469000000.0
$ws.Range("B32").Value = 1473000000.0

# Correct previously mis-entered values
$ws.Range("B27").Value = -612000000.0
$ws.Range("F14").Value = -1695000000.0
$ws.Range("F29").Value = -1722000000.0
